$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "QC-001-002"
$ws.Range("A16").Value = "QC-002-002"
$ws.Range("A21").Value = "QC-001-003"
$ws.Range("A22").Value = "QC-002-003"

$ws.Range("P27").Select()
